$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 736.1111
$ws.Range("I19").Value = 564.1667
$ws.Range("K19").Value = 564.1667
$ws.Range("M19").Value = -389.1667
$ws.Range("H43").Value = 460
$ws.Range("I43").Value = 460
$ws.Range("K43").Value = 460
$ws.Range("M43").Value = -391
$ws.Range("H116").Value = 2005
$ws.Range("I116").Value = 2005
$ws.Range("K116").Value = 2005
$ws.Range("M116").Value = 1437
$ws.Range("H132").Value = 3717.5
$ws.Range("I132").Value = 2951.25
$ws.Range("K132").Value = 8853.75
$ws.Range("M132").Value = -6323.75
$ws.Range("H135").Value = 785
$ws.Range("I135").Value = 785
$ws.Range("K135").Value = 7065
$ws.Range("M135").Value = -4530

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2086.7778
$ws.Range("I2").Value = 2130.3333
$ws.Range("K2").Value = 2130.3333
$ws.Range("M2").Value = -2017.3333
$ws.Range("H13").Value = 0
$ws.Range("J13").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("N13").Value = ""
$ws.Range("H45").Value = 5379.8
$ws.Range("I45").Value = 1700
$ws.Range("K45").Value = 1700
$ws.Range("M45").Value = -1323
$ws.Range("H61").Value = 7610.8887
$ws.Range("I61").Value = 4333
$ws.Range("K61").Value = 4333
$ws.Range("M61").Value = -4121
$ws.Range("H74").Value = 7125
$ws.Range("H77").Value = 7125
$ws.Range("H95").Value = 8441.4
$ws.Range("J95").Value = 8441.4
$ws.Range("L95").Value = 8441.4
$ws.Range("N95").Value = -13933.4
$ws.Range("H102").Value = 0
$ws.Range("I102").Value = 0
$ws.Range("K102").Value = 0
$ws.Range("M102").Value = ""
$ws.Range("H110").Value = 2558.75
$ws.Range("I110").Value = 2518.6365
$ws.Range("K110").Value = 2518.6365
$ws.Range("M110").Value = -473.6365000000001
$ws.Range("H116").Value = 2086.7778
$ws.Range("I116").Value = 2130.3333
$ws.Range("K116").Value = 2130.3333
$ws.Range("M116").Value = 163.6667000000002
$ws.Range("H122").Value = 3609.8
$ws.Range("I122").Value = 1699.6666
$ws.Range("K122").Value = 5098.9998
$ws.Range("M122").Value = -2648.9998
$ws.Range("H124").Value = 58490.11
$ws.Range("J124").Value = 58490.11
$ws.Range("L124").Value = 58490.11
$ws.Range("N124").Value = -68310.11
$ws.Range("H136").Value = 7610.8887
$ws.Range("I136").Value = 4333
$ws.Range("K136").Value = 12999
$ws.Range("M136").Value = -10449

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2086.7778
$ws.Range("I3").Value = 2130.3333
$ws.Range("K3").Value = 2130.3333
$ws.Range("M3").Value = -2016.3333
$ws.Range("H80").Value = 1598.5714
$ws.Range("I80").Value = 1488
$ws.Range("J80").Value = 1681.5
$ws.Range("K80").Value = 1488
$ws.Range("L80").Value = 1681.5
$ws.Range("M80").Value = -490
$ws.Range("N80").Value = -3677.5
$ws.Range("H83").Value = 1598.5714
$ws.Range("I83").Value = 1488
$ws.Range("J83").Value = 1681.5
$ws.Range("K83").Value = 7440
$ws.Range("L83").Value = 8407.5
$ws.Range("M83").Value = -2448
$ws.Range("N83").Value = -18391.5
$ws.Range("H99").Value = 1921.25
$ws.Range("I99").Value = 1921.25
$ws.Range("K99").Value = 1921.25
$ws.Range("M99").Value = -423.25
$ws.Range("H124").Value = 13997
$ws.Range("J124").Value = 13997
$ws.Range("L124").Value = 13997
$ws.Range("N124").Value = -23817

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H11").Value = 578
$ws.Range("I11").Value = 309
$ws.Range("J11").Value = 667.6667
$ws.Range("K11").Value = 309
$ws.Range("L11").Value = 667.6667
$ws.Range("M11").Value = -169
$ws.Range("N11").Value = -947.6667
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = ""
$ws.Range("N16").Value = ""
$ws.Range("H58").Value = 4100.231
$ws.Range("I58").Value = 2664
$ws.Range("K58").Value = 2664
$ws.Range("M58").Value = -2461
$ws.Range("H92").Value = 25989.834
$ws.Range("J92").Value = 25989.834
$ws.Range("L92").Value = 25989.834
$ws.Range("N92").Value = -30981.834
$ws.Range("H99").Value = 4748.75
$ws.Range("I99").Value = 4331.6665
$ws.Range("K99").Value = 4331.6665
$ws.Range("M99").Value = -2833.6665
$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = ""
$ws.Range("N113").Value = ""
$ws.Range("H126").Value = 4748.75
$ws.Range("I126").Value = 4331.6665
$ws.Range("K126").Value = 12994.9995
$ws.Range("M126").Value = -10524.9995
$ws.Range("H136").Value = 4100.231
$ws.Range("I136").Value = 2664
$ws.Range("K136").Value = 7992
$ws.Range("M136").Value = -5442

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 90914220
$ws.Range("I4").Value = 833.3333
$ws.Range("J4").Value = 125006750
$ws.Range("K4").Value = 2499.9999
$ws.Range("L4").Value = 375020250
$ws.Range("M4").Value = -2387.9999
$ws.Range("N4").Value = -375020474
$ws.Range("H34").Value = 1166.6666
$ws.Range("J34").Value = 1166.6666
$ws.Range("L34").Value = 3499.9998
$ws.Range("N34").Value = -3667.9998
$ws.Range("H39").Value = 100
$ws.Range("I39").Value = 100
$ws.Range("J39").Value = 0
$ws.Range("K39").Value = 300
$ws.Range("L39").Value = 0
$ws.Range("M39").Value = -6
$ws.Range("N39").Value = ""
$ws.Range("H122").Value = 428.2857
$ws.Range("I122").Value = 350
$ws.Range("J122").Value = 898
$ws.Range("K122").Value = 3150
$ws.Range("L122").Value = 8082
$ws.Range("M122").Value = -700
$ws.Range("N122").Value = -12982
$ws.Range("H124").Value = 4500
$ws.Range("J124").Value = 4500
$ws.Range("L124").Value = 13500
$ws.Range("N124").Value = -23320

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H52").Value = 40000
$ws.Range("J52").Value = 40000
$ws.Range("L52").Value = 40000
$ws.Range("N52").Value = -40518
$ws.Range("H101").Value = 29999.5
$ws.Range("J101").Value = 29999.5
$ws.Range("L101").Value = 29999.5
$ws.Range("N101").Value = -36489.5
$ws.Range("H102").Value = 5685.5713
$ws.Range("I102").Value = 5685.5713
$ws.Range("K102").Value = 5685.5713
$ws.Range("M102").Value = -4063.5713
$ws.Range("H123").Value = 39090.547
$ws.Range("J123").Value = 39090.547
$ws.Range("L123").Value = 39090.547
$ws.Range("N123").Value = -43990.547
$ws.Range("H132").Value = 9097.77
$ws.Range("I132").Value = 8355.916999999999
$ws.Range("K132").Value = 25067.751
$ws.Range("M132").Value = -22537.751
$ws.Range("H136").Value = 20000
$ws.Range("J136").Value = 20000
$ws.Range("L136").Value = 60000
$ws.Range("N136").Value = -65100

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H10").Value = 3248.25
$ws.Range("I10").Value = 5250
$ws.Range("K10").Value = 5250
$ws.Range("M10").Value = -5110
$ws.Range("H46").Value = 4300
$ws.Range("J46").Value = 4000
$ws.Range("L46").Value = 4000
$ws.Range("N46").Value = -4376

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H10").Value = 14186667
$ws.Range("I10").Value = 21255000
$ws.Range("K10").Value = 21255000
$ws.Range("M10").Value = -21254831
$ws.Range("H69").Value = 19998
$ws.Range("J69").Value = 19998
$ws.Range("L69").Value = 19998
$ws.Range("N69").Value = -21496
$ws.Range("H72").Value = 19998
$ws.Range("J72").Value = 19998
$ws.Range("L72").Value = 59994
$ws.Range("N72").Value = -67482
$ws.Range("H136").Value = 7474.625
$ws.Range("I136").Value = 1574.25
$ws.Range("K136").Value = 4722.75
$ws.Range("M136").Value = -2172.75
